$wb = $excel.ActiveWorkbook

# --- Sheet "Riscos" (sheet1) ---
$wsRiscos = $wb.Worksheets.Item("Riscos")

# Fill in Probabilidade (E3) and Impacto (F3) for the first risk row; the
# Severidade column (C3) recalculates automatically via its shared formula.
$wsRiscos.Range("E3").Value = 1
$wsRiscos.Range("F3").Value = 1

# Widen column C (Cod./Severidade) so the header text fits better.
$wsRiscos.Columns.Item(3).ColumnWidth = 11.175

# --- Sheet "Grafico" (sheet2) ---
$wsGrafico = $wb.Worksheets.Item("Grafico")
$wsGrafico.Activate()
$wsGrafico.Range("I10").Select()

# --- Make "Riscos" the active sheet/tab with H2 selected ---
$wsRiscos.Activate()
$wsRiscos.Range("H2").Select()
